$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of an existing data row (row 2) so the new rows pick up
# the same cell style (s="2") used throughout the table, without creating any
# new style entries.
$ws.Range("A2:G2").Copy()
$ws.Range("A34:G47").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New player rows (Liverpool squad) appended below the existing table.
$data = @(
    @("Steven Gerrard",   "Midfield", 34, "England", 8,  "Liverpool", 185),
    @("Coutinho",         "Midfield", 2,  "Brazil",  10, "Liverpool", 171),
    @("Jordan Henderson", "Midfield", 24, "England", 14, "Liverpool", 182),
    @("Adam Lallana",     "Midfield", 26, "England", 20, "Liverpool", 172),
    @("Lucas",            "Midfield", 27, "Brazil",  21, "Liverpool", 179),
    @("Emre Can",         "Midfield", 23, "Germany", 23, "Liverpool", 186),
    @("Joe Allen",        "Midfield", 24, "Wales",   24, "Liverpool", 168),
    @("Suso",             "Midfield", 20, "Spain",   30, "Liverpool", 176),
    @("Raheem Sterling",  "Midfield", 19, "England", 31, "Liverpool", 170),
    @("Lazar Marković",   "Midfield", 20, "Serbia",  50, "Liverpool", 174),
    @("Rickie Lambert",   "Forward",  32, "England", 9,  "Liverpool", 188),
    @("Daniel Sturridge", "Forward",  25, "England", 15, "Liverpool", 188),
    @("Fabio Borini",     "Forward",  23, "Italy",   29, "Liverpool", 180),
    @("Mario Balotelli",  "Forward",  24, "Italy",   45, "Liverpool", 189)
)

$row = 34
foreach ($player in $data) {
    $ws.Cells.Item($row, 1).Value = $player[0]
    $ws.Cells.Item($row, 2).Value = $player[1]
    $ws.Cells.Item($row, 3).Value = $player[2]
    $ws.Cells.Item($row, 4).Value = $player[3]
    $ws.Cells.Item($row, 5).Value = $player[4]
    $ws.Cells.Item($row, 6).Value = $player[5]
    $ws.Cells.Item($row, 7).Value = $player[6]
    $row = $row + 1
}

[void]$ws.Range("G38").Select()
